$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Race Time Calculation")
$ws.Range("G2").Value = "Laps"
$ws.Range("G2").Font.Bold = $true
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").Borders.Item(7).LineStyle = 1
Write-Host "done"
